$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 893.959796954314
$ws.Range("C2").Value = 914.327527352297
$ws.Range("D2").Value = 932.48099533437
$ws.Range("B3").Value = 542.403167060327
$ws.Range("C3").Value = 555.341696303482
$ws.Range("D3").Value = 565.733711191733
$ws.Range("B4").Value = 380.323851322205
$ws.Range("C4").Value = 388.551883228406
$ws.Range("D4").Value = 396.438908460383
$ws.Range("B5").Value = 562.786925261143
$ws.Range("C5").Value = 575.173533368636
$ws.Range("D5").Value = 587.235688442397
$ws.Range("B6").Value = 1463.13440410129
$ws.Range("C6").Value = 1496.51792386254
$ws.Range("D6").Value = 1524.10044018463
$ws.Range("B7").Value = 959.229275578443
$ws.Range("C7").Value = 978.439039936561
$ws.Range("D7").Value = 999.222130745871
$ws.Range("B8").Value = 982.59692269045
$ws.Range("C8").Value = 1005.76886292809
$ws.Range("D8").Value = 1027.27763959983
$ws.Range("B9").Value = 542.64043005521
$ws.Range("C9").Value = 553.818705434482
$ws.Range("D9").Value = 565.890718876197
$ws.Range("B10").Value = 490.372070171191
$ws.Range("C10").Value = 501.488853053238
$ws.Range("D10").Value = 511.245415130127
$ws.Range("B11").Value = 1300.50202332454
$ws.Range("C11").Value = 1330.58001923712
$ws.Range("D11").Value = 1356.86574691846
$ws.Range("B12").Value = 1688.66197001574
$ws.Range("C12").Value = 1721.9183512302
$ws.Range("D12").Value = 1761.11285553992
$ws.Range("B13").Value = 1652.7041145979
$ws.Range("C13").Value = 1687.12987953095
$ws.Range("D13").Value = 1722.7723375871
$ws.Range("B14").Value = 703.930195486456
$ws.Range("C14").Value = 718.142013370853
$ws.Range("D14").Value = 734.253716760245
$ws.Range("B15").Value = 1521.19605642188
$ws.Range("C15").Value = 1552.72205407014
$ws.Range("D15").Value = 1583.67729867611
$ws.Range("D16").Value = 409.102108546611
$ws.Range("D17").Value = 408.579504967266
$ws.Range("B18").Value = 335.548973561728
$ws.Range("C18").Value = 343.050697749874
$ws.Range("D18").Value = 349.807727657616
$ws.Range("B19").Value = 330.033598509987
$ws.Range("C19").Value = 337.300559275273
$ws.Range("D19").Value = 343.957866129146
$ws.Range("B20").Value = 1351.62385687935
$ws.Range("C20").Value = 1381.34132184368
$ws.Range("D20").Value = 1406.47388081024
$ws.Range("B21").Value = 1172.94407079287
$ws.Range("D22").Value = 1403.57875827806
$ws.Range("D23").Value = 1639.0916630286
$ws.Range("D24").Value = 1609.69584176948
$ws.Range("D25").Value = 2074.4287686286
$ws.Range("B26").Value = 1289.06757792852
$ws.Range("D27").Value = 1356.28077656945
$ws.Range("D28").Value = 2775.71596758281
$ws.Range("D29").Value = 2762.8788278928
$ws.Range("D30").Value = 2410.78076077037
$ws.Range("D31").Value = 2326.64847359856
$ws.Range("D32").Value = 544.899767759808
$ws.Range("D33").Value = 587.421578646044
$ws.Range("B34").Value = 606.460724426436
$ws.Range("C34").Value = 619.634988584914
$ws.Range("D34").Value = 632.870814245958
$ws.Range("B35").Value = 1270.86361964273
$ws.Range("C35").Value = 1297.18185671345
$ws.Range("D35").Value = 1326.65740952909
$ws.Range("B36").Value = 1025.27803543113
$ws.Range("C36").Value = 1049.2813938822
$ws.Range("D36").Value = 1069.61263844682
$ws.Range("B37").Value = 646.496802459432
$ws.Range("C37").Value = 660.989565224419
$ws.Range("D37").Value = 674.264994461587
$ws.Range("B38").Value = 1211.00362349335
$ws.Range("C38").Value = 1237.9607592873
$ws.Range("D38").Value = 1265.15888985961
$ws.Range("B39").Value = 1243.52730188362
$ws.Range("C39").Value = 1272.54989220209
$ws.Range("D39").Value = 1298.39582534197
$ws.Range("B40").Value = 1037.0011696789
$ws.Range("C40").Value = 1059.56321520753
$ws.Range("D40").Value = 1081.03664651992
$ws.Range("B41").Value = 1126.37189079231
$ws.Range("C41").Value = 1152.80951733676
$ws.Range("D41").Value = 1174.06105980764
$ws.Range("B42").Value = 1127.59882107495
$ws.Range("C42").Value = 1152.3254625559
$ws.Range("D42").Value = 1175.07133448128
$ws.Range("B43").Value = 461.400243767956
$ws.Range("C43").Value = 471.26786559647
$ws.Range("D43").Value = 481.054325227635
$ws.Range("B44").Value = 961.747245851213
$ws.Range("C44").Value = 984.610083240413
$ws.Range("D44").Value = 1003.45417026655
$ws.Range("B45").Value = 1478.52813454046
$ws.Range("D45").Value = 1543.82123305468
$ws.Range("B46").Value = 1472.05652025082
$ws.Range("C46").Value = 1507.10525267037
$ws.Range("D46").Value = 1536.32421387553
$ws.Range("B47").Value = 620.199299169799
$ws.Range("C47").Value = 634.104191653497
$ws.Range("D47").Value = 647.11476624933
$ws.Range("B48").Value = 523.824023356324
$ws.Range("C48").Value = 535.616089990898
$ws.Range("D48").Value = 546.16131553897
